$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.805.75"
$ws.Range("E2").Value = "  -0.98%  "

# Row 3
$ws.Range("D3").Value = "1.941.38"
$ws.Range("E3").Value = "  -0.95%  "

# Row 4
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.85"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  -2.07%  "

# Row 7
$style = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4884"
$ws.Range("D7").Style = $style
$ws.Range("E7").Value = "  -0.49%  "

# Row 8
$ws.Range("E8").Value = "  -0.73%  "

# Row 9
$style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06896"
$ws.Range("D9").Style = $style
$ws.Range("E9").Value = "  +0.91%  "

# Row 10
$style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.42"
$ws.Range("D10").Style = $style
$ws.Range("E10").Value = "  +1.75%  "

# Row 11
$style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "106.50"
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = "  -0.09%  "

# Row 12
$ws.Range("D12").Value = "1.939.27"
$ws.Range("E12").Value = "  +0.31%  "

# Row 13
$style = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07715"
$ws.Range("D13").Style = $style
$ws.Range("E13").Value = "  -0.47%  "

# Row 14
$style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.342"
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = "  -1.43%  "

# Row 15
$style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6987"
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = "  -1.73%  "

# Row 16
$style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "276.25"
$ws.Range("D16").Style = $style
$ws.Range("E16").Value = "  -3.27%  "

# Row 17
$ws.Range("D17").Value = "30.808.53"
$ws.Range("E17").Value = "  -1.04%  "

# Row 18
$style = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007719"
$ws.Range("D18").Style = $style
$ws.Range("E18").Value = "  -0.91%  "

# Row 19
$ws.Range("E19").Value = "  -1.07%  "

# Row 20
$ws.Range("D20").Value = "2.202.27"
$ws.Range("E20").Value = "  +0.63%  "

# Row 21
$ws.Range("E21").Value = "  -0.14%  "

# Row 22
$style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.450"
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = "  -1.97%  "

# Row 23
$style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = "  +0.03%  "

# Row 24
$style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.510"
$ws.Range("D24").Style = $style
$ws.Range("E24").Value = "  -1.35%  "

# Row 25
$style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.721"
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = "  -2.62%  "

# Row 26
$style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.96"
$ws.Range("D26").Style = $style
$ws.Range("E26").Value = "  -0.53%  "

# Row 27
$style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.63"
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = "  -2.10%  "

# Row 28
$style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.151"
$ws.Range("D28").Style = $style
$ws.Range("E28").Value = "  -2.21%  "

# Row 29
$style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1046"
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = "  -1.37%  "

# Row 30
$style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.386"
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = "  -3.91%  "

# Row 31
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.553"
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = "  -4.96%  "

# Row 32
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.550"
$ws.Range("D32").Style = $style
$ws.Range("E32").Value = "  -2.84%  "

# Row 33
$style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.366"
$ws.Range("D33").Style = $style
$ws.Range("E33").Value = "  -3.68%  "

# Row 34
$style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04845"
$ws.Range("D34").Style = $style
$ws.Range("E34").Value = "  -3.43%  "

# Row 35
$style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7505"
$ws.Range("D35").Style = $style
$ws.Range("E35").Value = "  -2.56%  "

# Row 36
$style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.157"
$ws.Range("D36").Style = $style
$ws.Range("E36").Value = "  -0.88%  "

# Row 37
$style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9996"
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = "  -0.04%  "

# Row 38
$style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.723"
$ws.Range("D38").Style = $style
$ws.Range("E38").Value = "  -0.39%  "

# Row 39
$style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01990"
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = "  -3.22%  "

# Row 41
$style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "78.56"
$ws.Range("D41").Style = $style
$ws.Range("E41").Value = "  +7.04%  "

# Row 42
$style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.469"
$ws.Range("D42").Style = $style
$ws.Range("E42").Value = "  +0.60%  "

# Row 43
$style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.094"
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = "  -2.03%  "

# Row 44
$style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9061"
$ws.Range("D44").Style = $style
$ws.Range("E44").Value = "  +2.28%  "

# Row 45
$style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "108.74"
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = "  -0.76%  "

# Row 46
$style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4394"
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = "  -1.68%  "

# Row 47
$style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9984"
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = "  -0.25%  "

# Row 48
$style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.747"
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = "  +3.46%  "

# Row 49
$style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "986.72"
$ws.Range("D49").Style = $style
$ws.Range("E49").Value = "  -3.72%  "

# Row 50
$style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1243"
$ws.Range("D50").Style = $style
$ws.Range("E50").Value = "  -2.25%  "

# Row 51
$style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.235"
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = "  -1.74%  "
